# Auto-generated edit script applying numeric updates per the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 10539.94
$ws.Range("I15").Value = 10539.94
$ws.Range("K15").Value = 31619.82
$ws.Range("M15").Value = -31450.82
$ws.Range("H17").Value = 1381653.1
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1381653.1
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4144959.3
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4145295.3
$ws.Range("H93").Value = 24527.4
$ws.Range("J93").Value = 24527.4
$ws.Range("L93").Value = 24527.4
$ws.Range("N93").Value = -29519.4
$ws.Range("H100").Value = 1985.9231
$ws.Range("I100").Value = 1701.25
$ws.Range("J100").Value = 2112.4443
$ws.Range("K100").Value = 1701.25
$ws.Range("L100").Value = 2112.4443
$ws.Range("M100").Value = -1160.25
$ws.Range("N100").Value = -3194.4443
$ws.Range("H108").Value = 60684
$ws.Range("J108").Value = 60684
$ws.Range("L108").Value = 60684
$ws.Range("N108").Value = -68364
$ws.Range("H112").Value = 1297.0392
$ws.Range("I112").Value = 620
$ws.Range("J112").Value = 1370.6305
$ws.Range("K112").Value = 1860
$ws.Range("L112").Value = 4111.8915
$ws.Range("M112").Value = -752
$ws.Range("N112").Value = -6327.8915
$ws.Range("H138").Value = 10418745
$ws.Range("I138").Value = 13890793
$ws.Range("J138").Value = 2600
$ws.Range("K138").Value = 41672379
$ws.Range("L138").Value = 7800
$ws.Range("M138").Value = -41667239
$ws.Range("N138").Value = -18080
$ws.Range("H139").Value = 275000
$ws.Range("J139").Value = 275000
$ws.Range("L139").Value = 275000
$ws.Range("N139").Value = -285280
$ws.Range("H141").Value = 1224.0233
$ws.Range("I141").Value = 1135.7
$ws.Range("J141").Value = 2401.6667
$ws.Range("K141").Value = 3407.1
$ws.Range("L141").Value = 7205.000100000001
$ws.Range("M141").Value = 1772.9
$ws.Range("N141").Value = -17565.0001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8599.322
$ws.Range("I32").Value = 8753.412
$ws.Range("J32").Value = 8047.8423
$ws.Range("K32").Value = 8753.412
$ws.Range("L32").Value = 8047.8423
$ws.Range("M32").Value = -8466.412
$ws.Range("N32").Value = -8621.8423
$ws.Range("H43").Value = 10564.111
$ws.Range("I43").Value = 11842
$ws.Range("J43").Value = 10404.375
$ws.Range("K43").Value = 11842
$ws.Range("L43").Value = 10404.375
$ws.Range("M43").Value = -11529
$ws.Range("N43").Value = -11030.375
$ws.Range("H61").Value = 8065747.5
$ws.Range("I61").Value = 9435201
$ws.Range("J61").Value = 1190.5555
$ws.Range("K61").Value = 9435201
$ws.Range("L61").Value = 1190.5555
$ws.Range("M61").Value = -9434989
$ws.Range("N61").Value = -1614.5555
$ws.Range("H74").Value = 6025582
$ws.Range("I74").Value = 7353987
$ws.Range("J74").Value = 3479.6
$ws.Range("K74").Value = 7353987
$ws.Range("L74").Value = 3479.6
$ws.Range("M74").Value = -7353113
$ws.Range("N74").Value = -5227.6
$ws.Range("H77").Value = 6025582
$ws.Range("I77").Value = 7353987
$ws.Range("J77").Value = 3479.6
$ws.Range("K77").Value = 36769935
$ws.Range("L77").Value = 17398
$ws.Range("M77").Value = -36765567
$ws.Range("N77").Value = -26134
$ws.Range("H107").Value = 21582.4
$ws.Range("J107").Value = 21582.4
$ws.Range("L107").Value = 21582.4
$ws.Range("N107").Value = -29262.4
$ws.Range("H122").Value = 5339.893
$ws.Range("I122").Value = 6600.25
$ws.Range("J122").Value = 2189
$ws.Range("K122").Value = 19800.75
$ws.Range("L122").Value = 6567
$ws.Range("M122").Value = -17350.75
$ws.Range("N122").Value = -11467
$ws.Range("H131").Value = 67435.836
$ws.Range("J131").Value = 67435.836
$ws.Range("L131").Value = 67435.836
$ws.Range("N131").Value = -77515.836
$ws.Range("H132").Value = 3522433
$ws.Range("I132").Value = 4311499
$ws.Range("J132").Value = 1984.4615
$ws.Range("K132").Value = 12934497
$ws.Range("L132").Value = 5953.3845
$ws.Range("M132").Value = -12931967
$ws.Range("N132").Value = -11013.3845
$ws.Range("H136").Value = 8065747.5
$ws.Range("I136").Value = 9435201
$ws.Range("J136").Value = 1190.5555
$ws.Range("K136").Value = 28305603
$ws.Range("L136").Value = 3571.6665
$ws.Range("M136").Value = -28303053
$ws.Range("N136").Value = -8671.666499999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2672.3333
$ws.Range("I105").Value = 1443.6818
$ws.Range("J105").Value = 3711.9614
$ws.Range("K105").Value = 1443.6818
$ws.Range("L105").Value = 3711.9614
$ws.Range("M105").Value = 303.3181999999999
$ws.Range("N105").Value = -7205.9614
$ws.Range("H134").Value = 2339.9656
$ws.Range("I134").Value = 1299.075
$ws.Range("J134").Value = 4653.0557
$ws.Range("K134").Value = 3897.225
$ws.Range("L134").Value = 13959.1671
$ws.Range("M134").Value = -1362.225
$ws.Range("N134").Value = -19029.1671
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 32969.6
$ws.Range("J96").Value = 32969.6
$ws.Range("L96").Value = 32969.6
$ws.Range("N96").Value = -38461.6
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H140").Value = 39520
$ws.Range("J140").Value = 39520
$ws.Range("L140").Value = 39520
$ws.Range("N140").Value = -49880
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 601.7241
$ws.Range("I5").Value = 259.0476
$ws.Range("J5").Value = 1501.25
$ws.Range("K5").Value = 777.1428
$ws.Range("L5").Value = 4503.75
$ws.Range("M5").Value = -665.1428
$ws.Range("N5").Value = -4727.75
$ws.Range("H48").Value = 1500
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1500
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 4500
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -5000
$ws.Range("H134").Value = 3934.946
$ws.Range("I134").Value = 1979.0476
$ws.Range("J134").Value = 6502.0625
$ws.Range("K134").Value = 5937.142800000001
$ws.Range("L134").Value = 19506.1875
$ws.Range("M134").Value = -867.1428000000005
$ws.Range("N134").Value = -29646.1875
$ws.Range("H135").Value = 601.7241
$ws.Range("I135").Value = 259.0476
$ws.Range("J135").Value = 1501.25
$ws.Range("K135").Value = 2331.4284
$ws.Range("L135").Value = 13511.25
$ws.Range("M135").Value = 203.5716000000002
$ws.Range("N135").Value = -18581.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2175.4583
$ws.Range("I100").Value = 2217
$ws.Range("J100").Value = 2150.5334
$ws.Range("K100").Value = 2217
$ws.Range("L100").Value = 2150.5334
$ws.Range("M100").Value = -1676
$ws.Range("N100").Value = -3232.5334
$ws.Range("H122").Value = 5534.2144
$ws.Range("I122").Value = 5182
$ws.Range("J122").Value = 6277.778
$ws.Range("K122").Value = 15546
$ws.Range("L122").Value = 18833.334
$ws.Range("M122").Value = -13096
$ws.Range("N122").Value = -23733.334
$ws.Range("H136").Value = 8335915.5
$ws.Range("I136").Value = 8929440
$ws.Range("J136").Value = 26576.25
$ws.Range("K136").Value = 26788320
$ws.Range("L136").Value = 79728.75
$ws.Range("M136").Value = -26785770
$ws.Range("N136").Value = -84828.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2343.611
$ws.Range("I100").Value = 2563.5
$ws.Range("J100").Value = 2068.75
$ws.Range("K100").Value = 5127
$ws.Range("L100").Value = 4137.5
$ws.Range("M100").Value = -4586
$ws.Range("N100").Value = -5219.5
$ws.Range("H136").Value = 765.3273
$ws.Range("I136").Value = 601.55316
$ws.Range("J136").Value = 1727.5
$ws.Range("K136").Value = 1804.65948
$ws.Range("L136").Value = 5182.5
$ws.Range("M136").Value = 745.3405199999997
$ws.Range("N136").Value = -10282.5
$ws.Range("H138").Value = 40205.8
$ws.Range("J138").Value = 44007.25
$ws.Range("L138").Value = 44007.25
$ws.Range("N138").Value = -54287.25
